$wb = $excel.ActiveWorkbook

# Sheet1: A1 value changes from 8 to 11
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A1").Value = 11

# Sheet2: Add rows 8, 9, 10
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("A8").Value = 2.5595299478784055
$ws2.Range("B8").Value = 2.8456292951976336
$ws2.Range("C8").Value = 1
$ws2.Range("D8").Value = 2.85674758928146
$ws2.Range("E8").Value = 3

$ws2.Range("A9").Value = 2.0497354976649786
$ws2.Range("B9").Value = 1.9006341561095608
$ws2.Range("C9").Value = 2
$ws2.Range("D9").Value = 1.7450227721776792
$ws2.Range("E9").Value = 3

$ws2.Range("A10").Value = 0.16120692144373794
$ws2.Range("B10").Value = 1.8619175241832844
$ws2.Range("C10").Value = 2
$ws2.Range("D10").Value = 1.8694675403716781
$ws2.Range("E10").Value = 3

$ws2.Range("A10:E10").Select() | Out-Null
$ws1.Select() | Out-Null


